$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StabilityCriteria")
$ws.Columns("A:F").AutoFit()
